$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.680.59'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '3.549.60'
$ws.Range('E3').Value = '  -1.74%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '198.94'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.73%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '587.51'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -2.97%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.615'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -1.85%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +0.67%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.631'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -3.09%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '52.34'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -3.05%  '
$ws.Range('E12').Value = '  -4.66%  '
$ws.Range('B13').Value = 'BitcoinCash'
$ws.Range('C13').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '691.92'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +16.68%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '9.36'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -2.05%  '
$ws.Range('D15').Value = '4.109.72'
$ws.Range('E15').Value = '  -1.94%  '
$ws.Range('D16').Value = '69.720.58'
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').Value = '3.557.87'
$ws.Range('E17').Value = '  -1.73%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '12.50'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -5.72%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '18.63'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -2.99%  '
$ws.Range('E20').Value = '  -0.67%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.973'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -2.16%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '18.07'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +1.50%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '108.64'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +5.77%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '5.21'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +0.60%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '4.43'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -4.37%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.96'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -3.04%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '10.30'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -4.27%  '
$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '9.74'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +1.22%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '33.76'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('B30').Value = 'dogwifhat'
$ws.Range('C30').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '4.40'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '6.94'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -2.69%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '11.96'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -3.13%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.112'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -3.90%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '62.26'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -1.55%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '3.810.04'
$ws.Range('E35').Value = '  -3.64%  '
$ws.Range('B36').Value = 'PEPE'
$ws.Range('C36').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D36').Value = '0.0₃0824'
$ws.Range('E36').Value = '  -8.19%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +0.03%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '3.70'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +4.16%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '501.71'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -5.27%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.95'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -6.67%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.375'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -4.65%  '
$ws.Range('E42').Value = '  +1.69%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '34.85'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -6.45%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.0461'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.05%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.95'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +2.95%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '3.37'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.138'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -2.14%  '
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.39%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '8.42'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -2.28%  '
$ws.Range('B50').Value = 'CoreDAO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.81'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +71.87%  '
$ws.Range('B51').Value = 'Jupiter'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qMgTxtv34+jupiter-jup'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '1.82'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +21.92%  '
